$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, [string]$text)
    if ($text -match '^-?\d+(\.\d+)?$') {
        $cellRange.Value = "'" + $text
        $cellRange.Style = "Normal"
    } else {
        $cellRange.Value = $text
    }
}

Set-TextValue $ws.Range("D2") "286.69"
Set-TextValue $ws.Range("D3") "21.19"
Set-TextValue $ws.Range("D4") "6.448"
Set-TextValue $ws.Range("D5") "0.06377"
Set-TextValue $ws.Range("D6") "3.599"
Set-TextValue $ws.Range("D8") "6.581"
Set-TextValue $ws.Range("D9") "0.8252"
Set-TextValue $ws.Range("D10") "0.01419"
Set-TextValue $ws.Range("D11") "0.1684"
Set-TextValue $ws.Range("D12") "0.08774"
Set-TextValue $ws.Range("D13") "0.03700"
Set-TextValue $ws.Range("D14") "0.03203"
Set-TextValue $ws.Range("D15") "0.09188"
Set-TextValue $ws.Range("D16") "3.712"
Set-TextValue $ws.Range("D17") "0.001659"
Set-TextValue $ws.Range("D18") "0.04746"
Set-TextValue $ws.Range("D19") "0.006221"
Set-TextValue $ws.Range("D20") "0.006292"
Set-TextValue $ws.Range("D21") "0.001072"
Set-TextValue $ws.Range("D22") "0.0001602"
Set-TextValue $ws.Range("D23") "3.786"
Set-TextValue $ws.Range("D24") "2.319"
Set-TextValue $ws.Range("D26") "0.1261"
Set-TextValue $ws.Range("D28") "0.0002709"
Set-TextValue $ws.Range("E28") "27UpBotsUBXT"
Set-TextValue $ws.Range("D40") "0.04795"
Set-TextValue $ws.Range("D41") "0.007159"
Set-TextValue $ws.Range("B42") "BKEXToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1116"
Set-TextValue $ws.Range("E42") "41BKEXTokenBKK"
Set-TextValue $ws.Range("B43") "CEJI"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003104"
Set-TextValue $ws.Range("E43") "42CEJICEJIWorstin24h"
Set-TextValue $ws.Range("D44") "0.01169"
Set-TextValue $ws.Range("D45") "0.00007101"
Set-TextValue $ws.Range("D47") "1.001"
Set-TextValue $ws.Range("D48") "0.004872"
Set-TextValue $ws.Range("D49") "0.00001502"
